# Automatic update of files.
# The three sightings currently on rows 4-6 get re-keyed: what used to be
# row 5's record is now row 4, what used to be row 6's record is now row
# 5, and what used to be row 4's record is now row 6 (sort-order column B
# is independently refreshed for every record, including row 7).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 (becomes the old row-5 record) ---------------------------------
$ws.Range("A4").Value = 111896689
$ws.Range("B4").Value = 90821
$ws.Range("D4").Value = "LC"
$ws.Range("E4").Value = 5964
$ws.Range("F4").Value = "Fjällig taggsvamp s.str."
$ws.Range("G4").Value = "Sarcodon imbricatus s.str."
$ws.Range("H4").Value = "(L.:Fr.) P.Karst."
$ws.Range("Q4").Value = 575759
# New empty placeholder cell that wasn't present before.
$ws.Range("I4").Copy($ws.Range("AF4"))

# --- Row 5 (becomes the old row-6 record) ---------------------------------
$ws.Range("A5").Value = 111896603
$ws.Range("B5").Value = 56575
$ws.Range("D5").Value = "NT"
$ws.Range("E5").Value = 103021
$ws.Range("F5").Value = "Talltita"
$ws.Range("G5").Value = "Poecile montanus"
$ws.Range("H5").Value = "(Conrad von Baldenstein, 1827)"
$ws.Range("Q5").Value = 575827
$ws.Range("R5").Value = 6703782
# New cells: K5/L5/N5 stay blank, M5 carries the activity note.
$ws.Range("I5").Copy($ws.Range("K5"))
$ws.Range("I5").Copy($ws.Range("L5"))
$ws.Range("M5").Value = "spel/sång"
$ws.Range("I5").Copy($ws.Range("N5"))
# AF5 no longer exists on this record.
$ws.Range("AF5").ClearContents()

# --- Row 6 (becomes the old row-4 record) ---------------------------------
$ws.Range("A6").Value = 111896686
$ws.Range("B6").Value = 89100
$ws.Range("E6").Value = 5754
$ws.Range("F6").Value = "Gultoppig fingersvamp"
$ws.Range("G6").Value = "Ramaria testaceoflava"
$ws.Range("H6").Value = "(Bres.) Corner"
$ws.Range("Q6").Value = 575755
$ws.Range("R6").Value = 6703742
# K6/L6/M6/N6 no longer exist on this record.
$ws.Range("K6").ClearContents()
$ws.Range("L6").ClearContents()
$ws.Range("M6").ClearContents()
$ws.Range("N6").ClearContents()

# --- Row 7: only the sort-order value changes ------------------------------
$ws.Range("B7").Value = 89924
